$wb = $excel.ActiveWorkbook

$ws2007 = $wb.Worksheets.Item("2007")
$ws2007.Activate()
$ws2007.Range("B2").Value = -2970.55
$ws2007.Range("B3").Select()

$wsNE = $wb.Worksheets.Item("Northern Europe")
$wsNE.Activate()
$wsNE.Range("B2").Value = -4120.55
$wsNE.Range("B3").Select()

$ws2013 = $wb.Worksheets.Item("2013")
$ws2013.Activate()
$ws2013.Range("B3").Value = -160.55000000000001
$ws2013.Range("B4").Select()
